$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.223.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.619.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.44%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.54%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '302.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3791'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.36%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08078'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.204'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.30%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.362'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.190'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001211'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.92%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.618.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.44%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06921'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.505'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.33%  '

# Row 21
$ws.Range("E21").Value = '  +0.27%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.37%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.296.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.555'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.53%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.071'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.85%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.91%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.252'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.59%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.58%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.801.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.062'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.51%  '

# Row 33
$ws.Range("B33").Value = 'WEMIXTOKEN'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.116'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.14%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.444'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.11%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '

# Row 36
$ws.Range("E36").Value = '  -2.66%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08676'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.49%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2466'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.87%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06883'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.44%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.833'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6871'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.50%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.309'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.41%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.03%  '

# Row 45
$ws.Range("E45").Value = '  +0.36%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6277'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.58%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.953'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.246'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.03%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07877'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.24%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.73%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.169'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.34%  '
